# Updates "想去人数" (F column) counts across all four sheets
# (展览, 演出, 本地生活, 全部类型) per "Update gh-pages to output generated at 34df19c".
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 278
$ws1.Range("F7").Value = 9332
$ws1.Range("F9").Value = 26
$ws1.Range("F11").Value = 1818
$ws1.Range("F12").Value = 39
$ws1.Range("F13").Value = 97
$ws1.Range("F14").Value = 2484
$ws1.Range("F16").Value = 3785
$ws1.Range("F17").Value = 276
$ws1.Range("F18").Value = 112
$ws1.Range("F20").Value = 193
$ws1.Range("F21").Value = 227
$ws1.Range("F22").Value = 187
$ws1.Range("F23").Value = 78
$ws1.Range("F24").Value = 52
$ws1.Range("F25").Value = 250
$ws1.Range("F26").Value = 507
$ws1.Range("F27").Value = 110
$ws1.Range("F28").Value = 1081
$ws1.Range("F29").Value = 444
$ws1.Range("F30").Value = 4296
$ws1.Range("F31").Value = 67
$ws1.Range("F32").Value = 75
$ws1.Range("F33").Value = 270
$ws1.Range("F34").Value = 61

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 35
$ws2.Range("F6").Value = 15

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 950

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 951
$ws4.Range("F6").Value = 35
$ws4.Range("F10").Value = 278
$ws4.Range("F11").Value = 9332
$ws4.Range("F13").Value = 26
$ws4.Range("F15").Value = 1818
$ws4.Range("F16").Value = 39
$ws4.Range("F17").Value = 97
$ws4.Range("F19").Value = 2484
$ws4.Range("F21").Value = 3785
$ws4.Range("F22").Value = 276
$ws4.Range("F23").Value = 112
$ws4.Range("F25").Value = 193
$ws4.Range("F26").Value = 227
$ws4.Range("F27").Value = 187
$ws4.Range("F29").Value = 78
$ws4.Range("F30").Value = 52
$ws4.Range("F31").Value = 250
$ws4.Range("F32").Value = 507
$ws4.Range("F33").Value = 110
$ws4.Range("F34").Value = 1081
$ws4.Range("F35").Value = 444
$ws4.Range("F36").Value = 4296
$ws4.Range("F37").Value = 67
$ws4.Range("F38").Value = 75
$ws4.Range("F39").Value = 270
$ws4.Range("F40").Value = 61
$ws4.Range("F41").Value = 15
